$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the CasesTab query text (cell B2): remove the stray leading space
# before "order By" on the final line, so the LIMIT 100 clause reads all
# tabs consistently (per commit message "Fixed limit-10 reading function
# to read all tabs").
$cell = $ws.Cells.Item(2, 2)
$old = $cell.Value2
$new = $old -replace "`n order By", "`norder By"
$cell.Value2 = $new

# Update the active selection to C2, matching the saved view state.
$ws.Range("C2").Select()
